# Add a new "2022" column to the renewable-energy-share table.
#
# The sheet has a year header row (row 4, D..R = 2007..2021) and two data
# rows below it (row 5: share %, row 6: hydropower production). We extend
# the series by one more year in column S, copying the formatting from the
# previous year's column (R) for each of the three rows and then writing
# in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 30

$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").Value = 11928.6

$ws.Application.CutCopyMode = $false

# The source workbook's selection moved to a single cell just past the new
# column rather than staying on the old R4:R6 range.
[void]$ws.Range("T3").Select()
